$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '36.587.56'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.009.64'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.20'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.73'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("E9").Value = '  +4.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.06'
$ws.Range("E10").Value = '  -2.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +7.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.104'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.885'
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.00'
$ws.Range("E14").Value = '  +12.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.25'
$ws.Range("E15").Value = '  -4.25%  '
$ws.Range("D16").Value = '2.303.21'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").Value = '2.010.60'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '36.501.09'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.86'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = '0.0₃0881'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.38'
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.20'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").Value = '  -8.03%  '
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("E27").Value = '  +5.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.42'
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.136'
$ws.Range("E29").Value = '  +23.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.14'
$ws.Range("E30").Value = '  +2.54%  '
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.02'
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.18'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  +4.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.51'
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("E36").Value = '  +9.25%  '
$ws.Range("E37").Value = '  -4.25%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("E40").Value = '  +15.13%  '
$ws.Range("E41").Value = '  -2.30%  '
$ws.Range("E42").Value = '  +3.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0218'
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.81'
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.57'
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.71'
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("D49").Value = '1.365.56'
$ws.Range("E49").Value = '  -4.42%  '
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").Value = '2.199.15'
$ws.Range("E51").Value = '  +0.18%  '
